# Update a set of numeric values in Sheet1 (result_data_RandomForest.xlsx)
# These correspond to updated algorithm results after renaming/updating the
# RandomForest algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.944900000000001
$ws.Range("C3").Value = -11.7648
$ws.Range("A12").Value = -21.60880000000001
$ws.Range("B14").Value = 6.615399999999997
$ws.Range("C20").Value = -11.66970000000001
$ws.Range("C25").Value = -13.3785
$ws.Range("B26").Value = 3.913300000000003
$ws.Range("A27").Value = -21.54009999999999
$ws.Range("C30").Value = -13.30799999999999
$ws.Range("B31").Value = 5.421900000000001
$ws.Range("A32").Value = -21.3762
$ws.Range("B35").Value = 9.345400000000005
$ws.Range("A36").Value = -19.6307
$ws.Range("B37").Value = 8.803600000000003
$ws.Range("A38").Value = -19.5406
$ws.Range("C44").Value = -13.19249999999999
$ws.Range("B45").Value = 6.265399999999997
$ws.Range("A46").Value = -21.4593
$ws.Range("C47").Value = -11.9326
$ws.Range("B52").Value = 5.028500000000002
$ws.Range("A54").Value = -21.56109999999999
$ws.Range("A55").Value = -22.1936
$ws.Range("A56").Value = -22.04460000000001
$ws.Range("B57").Value = 4.932899999999996
$ws.Range("C58").Value = -13.3065
$ws.Range("A67").Value = -21.53319999999998
$ws.Range("A69").Value = -21.67449999999998
$ws.Range("A72").Value = -21.914
$ws.Range("C78").Value = -11.10820000000001
$ws.Range("B81").Value = 6.574599999999999
$ws.Range("A83").Value = -21.6131
$ws.Range("B83").Value = 5.389900000000004
$ws.Range("C84").Value = -14.05359999999999
$ws.Range("A86").Value = -21.9938
$ws.Range("C89").Value = -10.9352
$ws.Range("A91").Value = -21.24550000000001
$ws.Range("C91").Value = -10.7346
$ws.Range("C92").Value = -11.07819999999999
$ws.Range("A93").Value = -21.24479999999999
$ws.Range("C96").Value = -13.2759
$ws.Range("A99").Value = -20.32029999999998
$ws.Range("B100").Value = 5.1115
$ws.Range("B102").Value = 8.237000000000002
$ws.Range("C102").Value = -13.5615
